# renames *UnitTest* to *ComponentTest*
#
# Targets the single slide in this deck (the "Test Driver" / "Component"
# diagram). Splits the "Ui Tests" and "Unit Tests" labels into two runs
# each (so "Tests" can keep its own run, matching how PowerPoint records
# an in-place text edit), nudges several shapes/connectors that were
# repositioned around the renamed box, and adds a new "JUnit" callout
# with its dashed connector arrow.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1. "Ui Tests" box (Rectangle 146) - split text into "Ui " + "Tests"
#    and shift it left.
# ---------------------------------------------------------------------
$rect146 = $s.Shapes.Item("Rectangle 146")
$rect146.TextFrame.TextRange.Text = "Ui Tests"
$rect146.TextFrame.TextRange.Characters(1, 3).Font.Bold = $true
$rect146.Left = 92.80897907795276

# ---------------------------------------------------------------------
# 2. "Ui Tests" box (Rectangle 147) - position only, text untouched.
# ---------------------------------------------------------------------
$rect147 = $s.Shapes.Item("Rectangle 147")
$rect147.Left = 96.28543477086615

# ---------------------------------------------------------------------
# 3. "Unit Tests" -> "Component Tests" box (Rectangle 148); reposition,
#    resize (taller, to fit the longer label) and re-split the text.
# ---------------------------------------------------------------------
$rect148 = $s.Shapes.Item("Rectangle 148")
$rect148.TextFrame.TextRange.Text = "Component Tests"
$rect148.TextFrame.TextRange.Characters(1, 10).Font.Bold = $true
$rect148.Left = 105.60818897637795
$rect148.Top = 189.31818897637794
$rect148.Height = 39.78488358976378

# ---------------------------------------------------------------------
# 4. Elbow Connector 164 (links Rectangle 146 to the Isosceles Triangle)
#    - re-routed: rotation/flip change plus new geometry.
# ---------------------------------------------------------------------
$conn164 = $s.Shapes.Item("Elbow Connector 164")
$conn164.Rotation = 90
$conn164.HorizontalFlip = -1
$conn164.VerticalFlip = -1
$conn164.Left = 136.35771943543307
$conn164.Top = 139.30071266141732
$conn164.Height = 12.814173228346457

# ---------------------------------------------------------------------
# 5. "AllUnitTestsSuite" box (Rectangle 169) - shifted left.
# ---------------------------------------------------------------------
$rect169 = $s.Shapes.Item("Rectangle 169")
$rect169.Left = 84.48622047244095

# ---------------------------------------------------------------------
# 6. Elbow Connector 175 (links Rectangle 169 to the Flowchart Decision)
#    - re-routed: rotation/flip change plus new geometry.
# ---------------------------------------------------------------------
$conn175 = $s.Shapes.Item("Elbow Connector 175")
$conn175.Rotation = 270
$conn175.HorizontalFlip = -1
$conn175.Left = 143.151970003937
$conn175.Top = 428.2332306464567
$conn175.Height = 0.4452755905511811

# ---------------------------------------------------------------------
# 7. Flowchart: Decision 279 - shifted left.
# ---------------------------------------------------------------------
$dec279 = $s.Shapes.Item("Flowchart: Decision 279")
$dec279.Left = 145.958031496063

# ---------------------------------------------------------------------
# 8. Elbow Connector 280 - shifted left to follow the decision shape.
# ---------------------------------------------------------------------
$conn280 = $s.Shapes.Item("Elbow Connector 280")
$conn280.Left = 140.4955978511811

# ---------------------------------------------------------------------
# 9. New "JUnit" textbox, cloned from the existing "TextBox 93" (Email
#    Server) callout so it inherits the identical style/lstStyle, then
#    moved into place and retexted.
# ---------------------------------------------------------------------
$jUnitBox = $s.Shapes.Item("TextBox 93").Duplicate().Item(1)
$jUnitBox.Name = "TextBox 57"
$jUnitBox.Left = 66.0
$jUnitBox.Top = 515.765594551181
$jUnitBox.TextFrame.TextRange.Text = "JUnit"

# ---------------------------------------------------------------------
# 10. New dashed arrow connector pointing at the JUnit textbox, cloned
#     from the matching "Straight Arrow Connector 302" so the line/
#     arrowhead formatting matches exactly.
# ---------------------------------------------------------------------
$jUnitArrow = $s.Shapes.Item("Straight Arrow Connector 302").Duplicate().Item(1)
$jUnitArrow.Name = "Straight Arrow Connector 58"
$jUnitArrow.Left = 120.0
$jUnitArrow.Top = 486.0
$jUnitArrow.Width = 0.00007874015748031496
$jUnitArrow.Height = 36.0
